$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "SUPER_ADMIN"
$ws.Range("B5").Value = "Super Admin"
$ws.Range("C5").Value = "Super Admin"
$ws.Range("D5").Value = "Admin"

$ws.Range("A5").Select()
